$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-06"

# Update the header label for the current-year column (I1) to match the new date
$ws.Range("I1").Value = "2022 (through 08-06)"

# Update the August figure for 2022 (row 9) and the running Total (row 14)
$ws.Range("I9").Value = 31
$ws.Range("I14").Value = 1001
